$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Remove the "source_id" survey row (row 5) entirely; this shifts all
# subsequent rows up by one, matching the target layout.
$ws.Rows.Item(5).Delete()
